# Apply updated values to "Means" and "Standard Deviations" sheets
# to include tri proximity tables.

$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# ---- Means sheet ----
$wsMeans.Range("E2").Value = 64
$wsMeans.Range("F2").Value = 60
$wsMeans.Range("G2").Value = 58

$wsMeans.Range("D3").Value = 16
$wsMeans.Range("F3").Value = 12
$wsMeans.Range("G3").Value = 11

$wsMeans.Range("D4").Value = 15
$wsMeans.Range("E4").Value = 24
$wsMeans.Range("G4").Value = 31

$wsMeans.Range("D5").Value = 26
$wsMeans.Range("E5").Value = 36
$wsMeans.Range("F5").Value = 39
$wsMeans.Range("G5").Value = 38

$wsMeans.Range("D6").Value = 70
$wsMeans.Range("E6").Value = 71
$wsMeans.Range("F6").Value = 69
$wsMeans.Range("G6").Value = 76

$wsMeans.Range("D7").Value = 5.9
$wsMeans.Range("E7").Value = 7.6
$wsMeans.Range("F7").Value = 8
$wsMeans.Range("G7").Value = 6.8

$wsMeans.Range("D8").Value = 6.8
$wsMeans.Range("E8").Value = 5.9
$wsMeans.Range("F8").Value = 6.1
$wsMeans.Range("G8").Value = 5

$wsMeans.Range("D9").Value = 46
$wsMeans.Range("E9").Value = 40
$wsMeans.Range("F9").Value = 38
$wsMeans.Range("G9").Value = 38

$wsMeans.Range("E10").Value = 0.45
$wsMeans.Range("G10").Value = 0.43

# ---- Standard Deviations sheet ----
$wsSD.Range("D2").Value = 28
$wsSD.Range("E2").Value = 23
$wsSD.Range("F2").Value = 23
$wsSD.Range("G2").Value = 24

$wsSD.Range("D3").Value = 30
$wsSD.Range("E3").Value = 25
$wsSD.Range("F3").Value = 21
$wsSD.Range("G3").Value = 18

$wsSD.Range("D4").Value = 14
$wsSD.Range("F4").Value = 19

$wsSD.Range("D5").Value = 24
$wsSD.Range("G5").Value = 28

$wsSD.Range("D6").Value = 29
$wsSD.Range("E6").Value = 32
$wsSD.Range("F6").Value = 31
$wsSD.Range("G6").Value = 30

$wsSD.Range("D7").Value = 7.2
$wsSD.Range("E7").Value = 10
$wsSD.Range("F7").Value = 9.8
$wsSD.Range("G7").Value = 8.7

$wsSD.Range("D8").Value = 12
$wsSD.Range("E8").Value = 8.5
$wsSD.Range("F8").Value = 7.7
$wsSD.Range("G8").Value = 6.7

$wsSD.Range("D9").Value = 26
$wsSD.Range("E9").Value = 21
$wsSD.Range("F9").Value = 18
$wsSD.Range("G9").Value = 15

$wsSD.Range("D10").Value = 0.11
$wsSD.Range("E10").Value = 0.088
$wsSD.Range("G10").Value = 0.082
